$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D36").Value = "Finding Optimal Augmentation"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/370"

$ws.Range("D37").Value = "[Paper Review] Vision Transformer with Deformable Attention"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=2288&mod=document&pageid=1"

$ws.Range("D50").Value = "functional model 사용법"
$ws.Range("E50").Value = "http://incredible.egloos.com/7544251"

$ws.Range("D51").Value = "[MariaDB] 컬럼 추가시 컬럼이 원하는 위치에 추가되게 하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/MariaDB-%EC%BB%AC%EB%9F%BC-%EC%B6%94%EA%B0%80%EC%8B%9C-%EC%BB%AC%EB%9F%BC%EC%9D%B4-%EC%9B%90%ED%95%98%EB%8A%94-%EC%9C%84%EC%B9%98%EC%97%90-%EC%B6%94%EA%B0%80%EB%90%98%EA%B2%8C-%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"
